$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 10917.3,            11579, 9364,  0.2545408407847087),
    @(1, 11066.46666666667,  11859, 9396,  0.2628584702809652),
    @(2, 10921.36666666667,  11558, 9991,  0.2602198918660482),
    @(3, 11558.53333333333,  12256, 10469, 0.2579678297042847),
    @(4, 10512.93333333333,  11319, 9611,  0.2498527447382609),
    @(5, 11466.9,            12369, 9922,  0.2650693655014038),
    @(6, 11013.23333333333,  11835, 10138, 0.2650287469228109),
    @(7, 10743.9,            11548, 10143, 0.2595669507980347),
    @(8, 10653.63333333333,  11415, 9151,  0.2526545921961467),
    @(9, 10609.66666666667,  11597, 9665,  0.2554513454437256)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
